$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 4..10 down to 5..11 (bottom-up) to make room for
# the new "Return to Supplier" report row, carrying values AND formatting
# without introducing new style entries (PasteSpecial(xlPasteAll) alone
# drops formatting in this engine, so we follow it with an explicit
# PasteSpecial(xlPasteFormats) from the same source).
for ($r = 10; $r -ge 4; $r--) {
    $src = $ws.Range("A" + $r + ":J" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":J" + ($r + 1))
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4104) | Out-Null
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# New row 4 keeps the same row-level formatting as the data rows below it;
# just overwrite the cell values/content with the new test case's data.
$ws.Range("A4").Value = "Pharmacy\Reports\Purchase\TC03ReturnToSupplierReport.py"
$ws.Range("B4").Value = "Norun"
$ws.Range("C4").Value = "PharmacyReport"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = "TC03"
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = "SNCH"
$ws.Range("I4:J4").ClearContents()

# Clean up stray I/J remnants picked up by the row-shifting copy/paste, then
# re-place the " " marker in column J: it originally sat on the row holding
# TC006 (old row 5), which is now row 6.
$ws.Range("I5:J11").ClearContents()
$ws.Range("J6").Value = " "

# Update the active/selected cell to match the edited workbook.
$ws.Range("G4").Select()
